$wb = $excel.ActiveWorkbook

$wsErros = $wb.Worksheets.Item("erros")

# Add the new rows/column to the "erros" sheet.
$wsErros.Range("C1").Value = "transacao"

$wsErros.Range("A3").Value = "DEV02"
$wsErros.Range("A4").Value = "DEV03"
$wsErros.Range("B3").Value = "`"Sem saldo para classificar/determinar`",`"Material/Saldo sem NFO`""
$wsErros.Range("C2").Value = "MR8M"
$wsErros.Range("B4").Value = "Id não retornada a zsgr_100 após estorno de pedido na ME22N"
$wsErros.Range("C3").Value = "ZSGR_100"
$wsErros.Range("C4").Value = "ZSGR_100"

# Resize column B to fit the newly-added (wider) content.
$wsErros.Columns.Item(2).EntireColumn.AutoFit()

# Make "erros" the active sheet/tab, and set the active cell selection.
$wsErros.Activate()
$wsErros.Range("C5").Select() | Out-Null
